$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E stay text so numeric-looking strings are not coerced to numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '31.454.68'
$ws.Range("E2").Value = '  +3.71%  '
$ws.Range("D3").Value = '2.013.05'
$ws.Range("E3").Value = '  +7.67%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '0.7661'
$ws.Range("E5").Value = '  +62.32%  '
$ws.Range("D6").Value = '259.16'
$ws.Range("E6").Value = '  +6.30%  '
$ws.Range("D7").Value = '0.9992'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '0.3608'
$ws.Range("E8").Value = '  +25.71%  '
$ws.Range("D9").Value = '28.66'
$ws.Range("E9").Value = '  +32.12%  '
$ws.Range("D10").Value = '0.07142'
$ws.Range("E10").Value = '  +10.27%  '
$ws.Range("D11").Value = '0.8565'
$ws.Range("E11").Value = '  +19.92%  '
$ws.Range("D12").Value = '0.08122'
$ws.Range("E12").Value = '  +4.24%  '
$ws.Range("D13").Value = '102.03'
$ws.Range("E13").Value = '  +6.06%  '
$ws.Range("D14").Value = '2.009.38'
$ws.Range("E14").Value = '  +7.43%  '
$ws.Range("D15").Value = '5.626'
$ws.Range("E15").Value = '  +9.86%  '
$ws.Range("D16").Value = '275.14'
$ws.Range("E16").Value = '  -3.22%  '
$ws.Range("D17").Value = '31.438.38'
$ws.Range("E17").Value = '  +3.69%  '
$ws.Range("D18").Value = '14.68'
$ws.Range("E18").Value = '  +13.05%  '
$ws.Range("D19").Value = '5.951'
$ws.Range("E19").Value = '  +13.20%  '
$ws.Range("D20").Value = '0.000008043'
$ws.Range("E20").Value = '  +7.54%  '
$ws.Range("D21").Value = '2.266.87'
$ws.Range("E21").Value = '  +7.37%  '
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '7.278'
$ws.Range("E24").Value = '  +16.56%  '
$ws.Range("D25").Value = '10.17'
$ws.Range("E25").Value = '  +12.97%  '
$ws.Range("D26").Value = '164.21'
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = '0.1467'
$ws.Range("E27").Value = '  +52.53%  '
$ws.Range("D28").Value = '20.19'
$ws.Range("E28").Value = '  +7.89%  '
$ws.Range("D29").Value = '2.388'
$ws.Range("E29").Value = '  +26.79%  '
$ws.Range("D30").Value = '1.624'
$ws.Range("E30").Value = '  +9.53%  '
$ws.Range("D31").Value = '4.664'
$ws.Range("E31").Value = '  +10.87%  '
$ws.Range("D32").Value = '1.355'
$ws.Range("E32").Value = '  +3.03%  '
$ws.Range("D33").Value = '4.423'
$ws.Range("E33").Value = '  +7.11%  '
$ws.Range("D34").Value = '0.05242'
$ws.Range("E34").Value = '  +8.84%  '
$ws.Range("D35").Value = '1.245'
$ws.Range("E35").Value = '  +11.42%  '
$ws.Range("D36").Value = '0.7682'
$ws.Range("E36").Value = '  +11.93%  '
$ws.Range("D37").Value = '2.791'
$ws.Range("E37").Value = '  +2.91%  '
$ws.Range("D38").Value = '0.02025'
$ws.Range("E38").Value = '  +7.36%  '
$ws.Range("D39").Value = '2.945'
$ws.Range("E39").Value = '  +3.71%  '
$ws.Range("D40").Value = '6.766'
$ws.Range("E40").Value = '  +8.62%  '
$ws.Range("D41").Value = '80.61'
$ws.Range("E41").Value = '  +6.68%  '
$ws.Range("D42").Value = '2.206'
$ws.Range("E42").Value = '  +14.75%  '
$ws.Range("D43").Value = '0.4815'
$ws.Range("E43").Value = '  +14.95%  '
$ws.Range("D44").Value = '0.8651'
$ws.Range("E44").Value = '  +5.01%  '
$ws.Range("D45").Value = '104.97'
$ws.Range("E45").Value = '  +4.26%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '7.736'
$ws.Range("E47").Value = '  +10.42%  '
$ws.Range("D48").Value = '9.988'
$ws.Range("E48").Value = '  +3.63%  '
$ws.Range("D49").Value = '0.4411'
$ws.Range("E49").Value = '  +13.34%  '
$ws.Range("D50").Value = '37.16'
$ws.Range("E50").Value = '  +5.92%  '

# Row 51: coin replaced (Maker -> Algorand)
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.1212'
$ws.Range("E51").Value = '  +17.03%  '
